$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# Columns A ("Date", e.g. "2023-05-29") and D ("Week", e.g. "22") hold
# text that *looks* like a date / a number. Assigning such literals
# straight to .Value makes Excel's automatic type detection coerce them
# into a date serial / a real number (and stamps a NumberFormat-driven
# style on the cell), which does not match how the rest of the sheet
# stores these columns (plain text, default style). Building the text
# via a TEXT() formula and then flattening it back to a static value
# with Copy + PasteSpecial (values only) keeps it as plain text without
# touching the cell's style.
$ws.Cells.Item($row, 1).Formula = '=TEXT("2023-05-29","@")'
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)

$ws.Cells.Item($row, 2).Value = "22:37:42"
$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).Formula = '=TEXT("22","@")'
$ws.Cells.Item($row, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4163)

$ws.Cells.Item($row, 5).Value = 119787
$ws.Cells.Item($row, 6).Value = 133521
$ws.Cells.Item($row, 7).Value = 158308
$ws.Cells.Item($row, 8).Value = 130536
$ws.Cells.Item($row, 9).Value = 174464
$ws.Cells.Item($row, 10).Value = 113825
$ws.Cells.Item($row, 11).Value = 198314
$ws.Cells.Item($row, 12).Value = 220243
$ws.Cells.Item($row, 13).Value = 172080
$ws.Cells.Item($row, 14).Value = 119833
$ws.Cells.Item($row, 15).Value = 38681
$ws.Cells.Item($row, 16).Value = 34857
$ws.Cells.Item($row, 17).Value = 50423
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36969
$ws.Cells.Item($row, 20).Value = -1
